# "Luftdruck in schwing korrigiert" - correct the measured air pressure
# value in the schwing1 worksheet. Every other changed cell (S2, V2:V11,
# W2:W11, X2:X11, Y2:Y11, Z2:Z11, AB2, AB5, AB8, ...) is a formula that
# depends - directly or indirectly - on R2, so fixing this single input
# cell and recalculating reproduces all of the derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schwing1")
$ws.Activate() | Out-Null

# Corrected air pressure reading (was 98000)
$ws.Range("R2").Value = 94888

# Make sure every dependent formula is recalculated with the corrected value
$excel.CalculateFull() | Out-Null

# Scroll the sheet view so column C becomes the left-most visible column
# instead of the previous selection on K2
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C1").Select() | Out-Null
